# Refresh the cryptos list price/volume columns (D = Price, E = Volume(1h))
# with the latest values, as produced by the scheduled scraper run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range("D2").Value = "68.122.44"
$ws.Range("E2").Value = "  -3.07%  "

# Row 3: Ethereum
$ws.Range("D3").Value = "3.807.29"
$ws.Range("E3").Value = "  +1.53%  "

# Row 4: TetherUSD
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.26%  "

# Row 5: BNB
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "595.37"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -3.77%  "

# Row 6: Solana
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "173.00"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -4.49%  "

# Row 7: LidoStakedEther
$ws.Range("D7").Value = "3.804.97"
$ws.Range("E7").Value = "  +1.57%  "

# Row 8: USDC
$ws.Range("E8").Value = "  -0.11%  "

# Row 9: XRP
$ws.Range("E9").Value = "  -1.28%  "

# Row 10: Dogecoin
$ws.Range("E10").Value = "  -3.82%  "

# Row 11: Toncoin
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.26"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.69%  "

# Row 12: Cardano
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.465"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -3.26%  "

# Row 13: Avalanche
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "38.06"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -5.04%  "

# Row 14: ShibaInu
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000245"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.80%  "

# Row 15: WrappedliquidstakedEther2.0
$ws.Range("D15").Value = "4.435.63"
$ws.Range("E15").Value = "  +1.14%  "

# Row 16: WrappedEther
$ws.Range("D16").Value = "3.797.17"
$ws.Range("E16").Value = "  +1.11%  "

# Row 17: WrappedBTC
$ws.Range("D17").Value = "68.157.96"
$ws.Range("E17").Value = "  -3.15%  "

# Row 18: TRON
$ws.Range("E18").Value = "  -4.48%  "

# Row 19: Polkadot
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.17"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -4.86%  "

# Row 20: Chainlink
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "16.09"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.60%  "

# Row 21: BitcoinCash
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "490.01"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.78%  "

# Row 22: Uniswap
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.30"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.11%  "

# Row 23: Polygon
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.735"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.82%  "

# Row 24: Litecoin
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "84.68"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.06%  "

# Row 25: Fetch.AI
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.38"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -8.85%  "

# Row 26: PEPE
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0000138"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +4.16%  "

# Row 27: InternetComputer(DFINITY)
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.22"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -6.16%  "

# Row 28: RenderToken
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.24"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -9.28%  "

# Row 29: Dai
$ws.Range("E29").Value = "  -0.09%  "

# Row 30: PancakeSwap
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.93"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.03%  "

# Row 31: ImmutableX
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.44"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.80%  "

# Row 32: EthereumClassic
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "32.84"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +7.52%  "

# Row 33: NEARProtocol
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "7.69"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.69%  "

# Row 34: Hedera
$ws.Range("E34").Value = "  -3.81%  "

# Row 35: FirstDigitalUSD
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.998"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.14%  "

# Row 36: Mantle
$ws.Range("E36").Value = "  -3.90%  "

# Row 37: Kaspa
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.137"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.87%  "

# Row 38: Filecoin
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.79"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -5.62%  "

# Row 39: TheGraph
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.327"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -7.36%  "

# Row 40: Bittensor
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "456.59"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +4.90%  "

# Row 41: OKB
$ws.Range("E41").Value = "  -2.14%  "

# Row 42: Stacks
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.00"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.43%  "

# Row 43: dogwifhat
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.90"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -5.32%  "

# Row 44: Cosmos
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.28"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.76%  "

# Row 45: Arweave
$ws.Range("E45").Value = "  -8.68%  "

# Row 46: Maker
$ws.Range("D46").Value = "2.830.90"
$ws.Range("E46").Value = "  -4.22%  "

# Row 47: Monero
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "139.69"
$ws.Range("D47").Style = "Normal"

# Row 48: USDe
$ws.Range("E48").Value = "  +0.03%  "

# Row 49: VeChain
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0352"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.99%  "

# Row 50: InjectiveProtocol
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "26.32"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -4.06%  "

# Row 51: EnergySwap
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "23.31"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +7.90%  "
